# The presentation currently uses the "Integral" theme (ppt/theme/theme2.xml,
# linked from the slide master) for all slides, while ppt/theme/theme1.xml
# (linked only from the notes master) holds the default Office Theme colours.
# The authored change swaps the two themes' colour schemes: the deck-facing
# theme becomes the plain "Office Theme" palette.
#
# PowerPoint's object model exposes the 12 DrawingML theme colour slots via
# Slide.ThemeColorScheme (Colors 1-12 = dk1, lt1, dk2, lt2, accent1-6, hlink,
# folHlink), which edits the theme part actually applied to the slide master
# (theme2.xml in this deck). We drive every slot to the target Office Theme
# RGB value using that API.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$tcs = $s.ThemeColorScheme

# Index order matches MsoThemeColorSchemeIndex:
# 1 dk1, 2 lt1, 3 dk2, 4 lt2, 5-10 accent1-6, 11 hlink, 12 folHlink
$targetHex = @(
    0x000000,   # 1  dk1
    0xFFFFFF,   # 2  lt1
    0x44546A,   # 3  dk2
    0xE7E6E6,   # 4  lt2
    0x5B9BD5,   # 5  accent1
    0xED7D31,   # 6  accent2
    0xA5A5A5,   # 7  accent3
    0xFFC000,   # 8  accent4
    0x4472C4,   # 9  accent5
    0x70AD47,   # 10 accent6
    0x0563C1,   # 11 hlink
    0x954F72    # 12 folHlink
)

for ($i = 0; $i -lt $targetHex.Length; $i++) {
    $idx = $i + 1
    $hex = $targetHex[$i]
    $r = ($hex -shr 16) -band 0xFF
    $g = ($hex -shr 8) -band 0xFF
    $b = $hex -band 0xFF
    $vbaRgb = $r + ($g * 256) + ($b * 65536)
    $tcs.Item($idx).RGB = $vbaRgb
}
